# Joung et al. Supp Data - shift sample-collection dates in column A
# back by exactly 1096 days (3 years, spanning a leap day) on both
# data sheets, leaving every other value/format untouched.

$wb = $excel.ActiveWorkbook

$dayShift = 1096

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [double]) {
            $cell.Value2 = $val - $dayShift
        }
    }
}
